# "Courbes finales dans excel"
# Slide 19 ("Broadcast Oriented Protocols"): the content placeholder's text
# currently renders with a dark-gray theme tint (schemeClr tx1 + lumMod 75% /
# lumOff 25%). Replace that with a plain solid black fill (srgbClr 000000)
# for every run in the placeholder, including the last run ("Relative
# Neighborhood Graph") which previously had no explicit fill at all.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(19)
$sh = $s.Shapes.Item(2)

$font = $sh.TextFrame.TextRange.Font
$font.Color.RGB = 0
